$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores every value as literal text in the source data
# (e.g. thousands-dot formatted "42.987.96" or trailing-zero "310.00"), so
# force text format on any updated Price cell whose new value would
# otherwise be auto-detected as a number, to avoid losing formatting/precision.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '42.987.96'
$ws.Range("E2").Value = '  +2.21%  '
$ws.Range("D3").Value = '2.304.03'
$ws.Range("E3").Value = '  +1.74%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '310.00'
$ws.Range("E5").Value = '  +1.53%  '
$ws.Range("D6").Value = '100.35'
$ws.Range("E6").Value = '  +4.82%  '
$ws.Range("D7").Value = '0.535'
$ws.Range("E7").Value = '  +1.42%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").Value = '0.517'
$ws.Range("E9").Value = '  +5.63%  '
$ws.Range("D10").Value = '36.05'
$ws.Range("E10").Value = '  +2.60%  '
$ws.Range("D11").Value = '0.0821'
$ws.Range("E11").Value = '  +3.94%  '
$ws.Range("E12").Value = '  +0.62%  '
$ws.Range("D13").Value = '7.16'
$ws.Range("E13").Value = '  +7.90%  '
$ws.Range("D14").Value = '2.665.29'
$ws.Range("E14").Value = '  +2.05%  '
$ws.Range("D15").Value = '14.91'
$ws.Range("E15").Value = '  +3.68%  '
$ws.Range("D16").Value = '2.320.79'
$ws.Range("E16").Value = '  +1.51%  '
$ws.Range("D17").Value = '0.802'
$ws.Range("E17").Value = '  +1.22%  '
$ws.Range("D18").Value = '42.952.21'
$ws.Range("E18").Value = '  +2.39%  '
$ws.Range("D19").Value = '12.53'
$ws.Range("E19").Value = '  +0.94%  '
$ws.Range("D20").Value = '0.0₃0920'
$ws.Range("E20").Value = '  +1.94%  '
$ws.Range("D21").Value = '6.07'
$ws.Range("E21").Value = '  +1.73%  '
$ws.Range("D22").Value = '68.22'
$ws.Range("E22").Value = '  +0.81%  '
$ws.Range("D23").Value = '239.76'
$ws.Range("E23").Value = '  +0.78%  '
$ws.Range("E24").Value = '  +4.67%  '
$ws.Range("D25").Value = '2.61'
$ws.Range("E25").Value = '  +1.21%  '
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").Value = '24.14'
$ws.Range("E27").Value = '  +1.92%  '
$ws.Range("D28").Value = '38.72'
$ws.Range("E28").Value = '  +5.47%  '
$ws.Range("D29").Value = '9.65'
$ws.Range("E29").Value = '  +1.59%  '
$ws.Range("D30").Value = '2.12'
$ws.Range("E30").Value = '  +0.35%  '
$ws.Range("D31").Value = '168.58'
$ws.Range("E31").Value = '  +5.88%  '
$ws.Range("D32").Value = '5.34'
$ws.Range("E32").Value = '  +2.16%  '
$ws.Range("E33").Value = '  +0.10%  '
$ws.Range("D34").Value = '3.13'
$ws.Range("E34").Value = '  -1.68%  '
$ws.Range("D35").Value = '17.69'
$ws.Range("E35").Value = '  +3.76%  '
$ws.Range("E36").Value = '  +0.21%  '
$ws.Range("E37").Value = '  +0.33%  '
$ws.Range("E38").Value = '  +0.47%  '
$ws.Range("E39").Value = '  +0.47%  '
$ws.Range("E40").Value = '  +1.55%  '
$ws.Range("D41").Value = '4.23'
$ws.Range("E41").Value = '  +5.19%  '
$ws.Range("E42").Value = '  -4.93%  '
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").Value = '0.0289'
$ws.Range("E43").Value = '  +2.07%  '
$ws.Range("D44").Value = '1.967.67'
$ws.Range("E44").Value = '  -0.75%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '19.21'
$ws.Range("E45").Value = '  +0.94%  '
$ws.Range("D46").Value = '3.02'
$ws.Range("E46").Value = '  +3.12%  '
$ws.Range("D47").Value = '9.79'
$ws.Range("E47").Value = '  -1.28%  '
$ws.Range("D48").Value = '2.97'
$ws.Range("E48").Value = '  +18.40%  '
$ws.Range("D49").Value = '55.12'
$ws.Range("E49").Value = '  +3.65%  '
$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D50").Value = '1.55'
$ws.Range("E50").Value = '  +2.52%  '
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").Value = '2.531.69'
$ws.Range("E51").Value = '  +1.80%  '
